# Status Update.xlsx - add the September 2020 status block and drop the
# old trailing placeholder rows (48-55), matching the "All Services
# Implemented with Symmetric Key Authorization" weekly status entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Remove the old trailing rows (48-55) entirely --------------------
$ws.Range("A48:F55").EntireRow.Delete()

# --- 2. New header row (row 49), formatted like the other section headers
$ws.Range("A40:F40").Copy()
$ws.Range("A49:F49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A49").Value = "Date"
$ws.Range("B49").Value = "Saurabh Mehra"
$ws.Range("C49").Value = "Mansi Verma"
$ws.Range("D49").Value = "Vikas Pandey"
$ws.Range("E49").Value = "Shivam Singh"
$ws.Range("F49").Value = "Overall Status"

# --- 3. New data row (row 51): date + status text ------------------------
$ws.Range("A42").Copy()
$ws.Range("A51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A51").Value = (New-Object DateTime(2020, 9, 11))

$ws.Range("B51").Value = "working On Cloud Deployment"
$ws.Range("C51").Value = "working On Cloud Deployment"
$ws.Range("D51").Value = "working On Cloud Deployment"
$ws.Range("E51").Value = "working On Cloud Deployment"
$ws.Range("F51").Value = "Design and Coding is Completed"

# --- 4. Rows 52-53: continuing status notes in column F ------------------
$ws.Range("F52").Value = "Nunit Testing Completed"
$ws.Range("F53").Value = "working On Cloud Deployment"

# --- 5. Row 55: leftover formatted placeholder cell -----------------------
$ws.Range("B50").Copy()
$ws.Range("B55").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Rows.Item(55).RowHeight = 15.75

$excel.CutCopyMode = $false

# --- 6. Update the sheet view to match the new scroll/selection position -
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("B57").Select()
